# Update the localization status report: the files
#   07b2a4e1-22d5-4332-9f70-73ce01e95d6e.md
#   311b0043-e1a3-4fa5-bb89-00f5c2e4373a.md
# have moved from "Ready for handoff" to "In Translation" for both the
# zh-cn and de-de locales. Update the Overview sheet (zh-cn/de-de status
# columns) as well as each locale's own detail sheet (Status column).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
# Row 3 = 07b2a4e1-22d5-4332-9f70-73ce01e95d6e.md, Row 4 = 311b0043-e1a3-4fa5-bb89-00f5c2e4373a.md
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
